# Update "想去人数" (interest count) values for two events that appear on
# multiple sheets: "南宁·草莓动漫节" (995 -> 1002) and
# "南宁·第一届ANE·DACG动漫嘉年华" (551 -> 554).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): row 4 = 草莓动漫节, row 5 = ANE·DACG动漫嘉年华
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1002
$ws1.Range("F5").Value = 554

# Sheet "全部类型" (all types): row 4 = 草莓动漫节, row 6 = ANE·DACG动漫嘉年华
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1002
$ws4.Range("F6").Value = 554
